$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new label columns (I, J) mirroring existing categories "E" and "F"
# (part of renaming a duplicate heading for the 2nd color group)
$ws.Range("I1").Value = "E"
$ws.Range("J1").Value = "F"

# Row 2: add extra tallies in the new/other columns
$ws.Range("C2").Value = 1
$ws.Range("H2").Value = 1
$ws.Range("J2").Value = 1

# Row 4: add tally in new column I
$ws.Range("I4").Value = 1

# Row 5: add tallies in columns C and J
$ws.Range("C5").Value = 1
$ws.Range("J5").Value = 1

# Row 6: add tally in column C
$ws.Range("C6").Value = 1

# Row 7: add tally in new column J
$ws.Range("J7").Value = 1

# Update the selected cell to match the saved view state
$ws.Range("K5").Select()
